$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 21) below the existing data (last existing row is 20),
# mirroring the structure/style of the preceding rows.
$row = 21

$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value = "Maule"

# Column D holds a date serial value; copy the style from the cell above so the
# date number format (and other style attributes) stays consistent.
$ws.Cells.Item($row, 4).Value = 44585
$ws.Cells.Item($row - 1, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 4).Value = 44585

$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 200
$ws.Cells.Item($row, 11).Value = 12000
$ws.Cells.Item($row, 12).Value = 12000
$ws.Cells.Item($row, 13).Value = 12000
$ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 667
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"

$wb.Save()
